# Add three new localization entries (rows 32-34) to Sheet1:
#   MINIGAME-LIMIT, NOT-ENOUGH, NOT-NULL
# These back the "coins-limit", "list items" and "login before play"
# features referenced in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 - minigame coin limit message
$ws.Range("A32").Value = "MINIGAME-LIMIT"
$ws.Range("B32").Value = "YOU CAN ONLY PLAY WHEN YOUR COIN IS UNDER 30!"
$ws.Range("C32").Value = "BẠN CHỈ ĐƯỢC CHƠI KHI CÓ ÍT HƠN 30 COIN!"

# Row 33 - not enough coins message
$ws.Range("A33").Value = "NOT-ENOUGH"
$ws.Range("B33").Value = "NOT ENOUGH COINS"
$ws.Range("C33").Value = "KHÔNG ĐỦ XU ĐỂ CHƠI"

# Row 34 - stake must not be empty message
$ws.Range("A34").Value = "NOT-NULL"
$ws.Range("B34").Value = "YOU MUST FILL YOUR STAKE!"
$ws.Range("C34").Value = "BẠN PHẢI NHẬP SỐ TIỀN CƯỢC!"

# Mirror the author's final selection/view state on save.
$ws.Range("C34").Select()
